# Updated symbol list on Tue Dec 13 03:47:45 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) quotes for most coins, and fix the mixed-up
# CEJI / BKEXToken rows (42 and 43) which had swapped data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates -------------------------------------------
# Column D is stored as text, so numeric-looking values are entered with a
# leading apostrophe to keep them as text instead of being converted to
# numbers.
$priceUpdates = @{
    2  = "269.15"
    3  = "21.11"
    4  = "6.254"
    5  = "0.06202"
    6  = "3.564"
    7  = "6.535"
    8  = "1.433"
    9  = "0.8247"
    10 = "0.1655"
    11 = "0.08241"
    12 = "0.03556"
    14 = "0.09191"
    15 = "3.775"
    16 = "0.001629"
    17 = "0.04681"
    18 = "0.006327"
    19 = "0.006190"
    20 = "0.001069"
    22 = "3.724"
    23 = "2.256"
    28 = "0.0002715"
    40 = "0.04695"
    41 = "0.007017"
    44 = "0.01049"
    45 = "0.00006220"
    46 = "0.0009906"
    47 = "0.00000000751"
    48 = "0.9911"
    50 = "0.00001901"
    51 = "0.01241"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# --- Rows 42/43: CEJI and BKEXToken were swapped -------------------------
# Row 42 becomes BKEXToken, row 43 becomes CEJI.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1120"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003521"
$ws.Range("E43").Value = "42CEJICEJI"
